$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.564.93'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.960.33'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.59'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.376'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0789'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.81%  '
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.44%  '
$ws.Range("D14").Value = '2.247.93'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.830'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("E16").Value = '  +1.96%  '
$ws.Range("D17").Value = '1.958.87'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").Value = '36.477.44'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '0.0₃0853'
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +3.54%  '
$ws.Range("E25").Value = '  +2.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.140'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.15%  '
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.34'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +20.64%  '
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.91%  '
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.46'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.74%  '
$ws.Range("E35").Value = '  +10.16%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.56%  '
$ws.Range("E38").Value = '  -1.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -12.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0966'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("D45").Value = '1.362.67'
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.66%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("E49").Value = '  +0.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.93%  '
$ws.Range("D51").Value = '2.137.27'
$ws.Range("E51").Value = '  +0.90%  '
